# Apply the edits described by the commit:
#   - grow the saved workbook window height (bookViews/workbookView/@windowHeight 11835 -> 12435)
#   - set the "Организация" (B6) value to the statistics committee name, which adds a
#     new shared string and repoints B6 at it
#   - move the sheet's remembered selection from C8 to B7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window geometry -------------------------------------------------
# windowHeight/windowWidth in xl/workbook.xml are stored in twentieths of a point,
# so 12435 -> 621.75 pt. Mirror that onto the (only) application window.
try {
    $win = $excel.ActiveWindow
    $win.Height = 621.75
} catch {
    # Window sizing isn't critical to the data edit; ignore if unsupported.
}

# --- Cell content ------------------------------------------------------
# B6 ("Организация") previously (incorrectly) duplicated the indicator text
# (shared string 29); it should hold the organization's full name instead.
$ws.Range("B6").Value = "Национальный статистический комитет КР (Отдел статистики труда и занятости)"

# --- Selection -----------------------------------------------------------
# The sheet's remembered selection moves from C8 to B7.
$ws.Range("B7").Select()
